$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins / Losses / Ties added to columns AD, AE, AF
# Copy header styling (bold font, thin border, centered alignment) from the
# adjacent existing header cell (AC1) so the new headers match the rest of
# row 1, then set the header text.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill team record (Wins=86, Losses=76, Ties=0) for every data row (2-42)
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 86   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 76   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
